# Apply the edits described by the diff to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for two additional data rows (the table grows from 6 to 8 rows) ---
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# --- Header row ---
$ws.Range("B1").Value = "Variable"
$ws.Range("C1").Value = "Detalle"
$ws.Range("D1").Value = "Fuente"
$ws.Range("E1").Value = "dates"

# --- Data rows (B:E), fully restated to match the new dataset layout ---
$ws.Range("B2").Value = " cpi"
$ws.Range("C2").Value = "IPC de Guatemala"
$ws.Range("D2").Value = "INE"
$ws.Range("E2").Value = "01/01/1995 -01/04/2022"

$ws.Range("B3").Value = "cpi_usa"
$ws.Range("C3").Value = "CPI for All Urban Consumers (CPI-U)"
$ws.Range("D3").Value = "U.S. Bureau Of Labor Statistics"
$ws.Range("E3").Value = "01/01/1995 -01/04/2022"

$ws.Range("B4").Value = "i_pm"
$ws.Range("C4").Value = "Tasa de Interés Líder de Polítca Monetaria"
$ws.Range("D4").Value = "Banco de Guatemala"
$ws.Range("E4").Value = "01/01/1995 -01/04/2022"

$ws.Range("B5").Value = "ner_gt"
$ws.Range("C5").Value = "Tipo de Cambio Nominal"
$ws.Range("D5").Value = "Banco de Guatemala"
$ws.Range("E5").Value = "01/01/1995 -01/04/2022"

$ws.Range("B6").Value = "pce_us"
$ws.Range("C6").Value = "PCE Inflation USA"
$ws.Range("D6").Value = "U.S. Bureau Of Labor Statistics"
$ws.Range("E6").Value = "01/01/1995 - 01/03/2022"

$ws.Range("B7").Value = "pce_us_core"
$ws.Range("C7").Value = "PCE core inflation USA"
$ws.Range("D7").Value = "U.S. Bureau Of Labor Statistics"
$ws.Range("E7").Value = "01/01/1995 - 01/03/2022"

$ws.Range("B8").Value = "imae_tc"
$ws.Range("C8").Value = "IMAE gt"
$ws.Range("D8").Value = "Banco de Guatemala"
$ws.Range("E8").Value = "01/01/2001 - 01/04/2022"

# --- Column E width to fit its content (bestFit) ---
$ws.Columns.Item(5).ColumnWidth = 21

# --- Selection matches the post-edit state recorded in the workbook view ---
$ws.Range("E3:E5").Select()
